$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1726.33
$ws.Range("I15").Value = 1726.33
$ws.Range("K15").Value = 5178.99
$ws.Range("M15").Value = -5009.99
$ws.Range("H76").Value = 2318009.8
$ws.Range("I76").Value = 3315.5557
$ws.Range("K76").Value = 3315.5557
$ws.Range("M76").Value = -3000.5557
$ws.Range("H79").Value = 2318009.8
$ws.Range("I79").Value = 3315.5557
$ws.Range("K79").Value = 3315.5557
$ws.Range("M79").Value = -2223.5557
$ws.Range("H98").Value = 1000.625
$ws.Range("I98").Value = 1000.625
$ws.Range("K98").Value = 1000.625
$ws.Range("M98").Value = 497.375
$ws.Range("H101").Value = 232.5
$ws.Range("J101").Value = 200
$ws.Range("L101").Value = 600
$ws.Range("N101").Value = -3844
$ws.Range("H113").Value = 142862540
$ws.Range("I113").Value = 250002190
$ws.Range("K113").Value = 250002190
$ws.Range("M113").Value = -249998936
$ws.Range("H122").Value = 1000.625
$ws.Range("I122").Value = 1000.625
$ws.Range("K122").Value = 3001.875
$ws.Range("M122").Value = -551.875
$ws.Range("H129").Value = 334019
$ws.Range("J129").Value = 417417.5
$ws.Range("L129").Value = 1252252.5
$ws.Range("N129").Value = -1262252.5
$ws.Range("H132").Value = 3645.25
$ws.Range("I132").Value = 4019.2083
$ws.Range("J132").Value = 1401.5
$ws.Range("K132").Value = 12057.6249
$ws.Range("L132").Value = 4204.5
$ws.Range("M132").Value = -9527.624899999999
$ws.Range("N132").Value = -9264.5
$ws.Range("H135").Value = 33341476
$ws.Range("I135").Value = 1690.5
$ws.Range("J135").Value = 100021050
$ws.Range("K135").Value = 15214.5
$ws.Range("L135").Value = 900189450
$ws.Range("M135").Value = -12679.5
$ws.Range("N135").Value = -900194520

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1128.5
$ws.Range("I2").Value = 1216.2
$ws.Range("J2").Value = 690
$ws.Range("K2").Value = 1216.2
$ws.Range("L2").Value = 690
$ws.Range("M2").Value = -1103.2
$ws.Range("N2").Value = -916
$ws.Range("H32").Value = 7540.2617
$ws.Range("I32").Value = 6361.1333
$ws.Range("J32").Value = 10193.3
$ws.Range("K32").Value = 6361.1333
$ws.Range("L32").Value = 10193.3
$ws.Range("M32").Value = -6074.1333
$ws.Range("N32").Value = -10767.3
$ws.Range("H112").Value = 27727.715
$ws.Range("J112").Value = 27727.715
$ws.Range("L112").Value = 27727.715
$ws.Range("N112").Value = -30681.715
$ws.Range("H116").Value = 1128.5
$ws.Range("I116").Value = 1216.2
$ws.Range("J116").Value = 690
$ws.Range("K116").Value = 1216.2
$ws.Range("L116").Value = 690
$ws.Range("M116").Value = 1077.8
$ws.Range("N116").Value = -5278

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1128.5
$ws.Range("I3").Value = 1216.2
$ws.Range("J3").Value = 690
$ws.Range("K3").Value = 1216.2
$ws.Range("L3").Value = 690
$ws.Range("M3").Value = -1102.2
$ws.Range("N3").Value = -918
$ws.Range("H94").Value = 567.3871
$ws.Range("J94").Value = 977.5
$ws.Range("L94").Value = 977.5
$ws.Range("N94").Value = -1879.5
$ws.Range("H107").Value = 1338.0769
$ws.Range("I107").Value = 1084.8334
$ws.Range("K107").Value = 1084.8334
$ws.Range("M107").Value = 835.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2052.75
$ws.Range("I16").Value = 1605.5
$ws.Range("K16").Value = 1605.5
$ws.Range("M16").Value = -1318.5
$ws.Range("H31").Value = 4310.115
$ws.Range("J31").Value = 5782.9443
$ws.Range("L31").Value = 5782.9443
$ws.Range("N31").Value = -6372.9443
$ws.Range("H34").Value = 4310.115
$ws.Range("J34").Value = 5782.9443
$ws.Range("L34").Value = 5782.9443
$ws.Range("N34").Value = -6186.9443
$ws.Range("H76").Value = 4450
$ws.Range("I76").Value = 4450
$ws.Range("K76").Value = 4450
$ws.Range("M76").Value = -4135
$ws.Range("H79").Value = 4450
$ws.Range("I79").Value = 4450
$ws.Range("K79").Value = 4450
$ws.Range("M79").Value = -3358
$ws.Range("H107").Value = 1347.7916
$ws.Range("I107").Value = 766.3077
$ws.Range("K107").Value = 766.3077
$ws.Range("M107").Value = 1153.6923
$ws.Range("H113").Value = 2052.75
$ws.Range("I113").Value = 1605.5
$ws.Range("K113").Value = 1605.5
$ws.Range("M113").Value = 564.5
$ws.Range("N113").Value = -6840
$ws.Range("H132").Value = 2211.4194
$ws.Range("I132").Value = 1288.1818
$ws.Range("K132").Value = 3864.5454
$ws.Range("M132").Value = -1334.5454
$ws.Range("H134").Value = 1120.2413
$ws.Range("I134").Value = 1003.625
$ws.Range("J134").Value = 1680
$ws.Range("K134").Value = 3010.875
$ws.Range("L134").Value = 5040
$ws.Range("M134").Value = -475.875
$ws.Range("N134").Value = -10110

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1196.3334
$ws.Range("I122").Value = 328.6
$ws.Range("J122").Value = 1424.6842
$ws.Range("K122").Value = 2957.4
$ws.Range("L122").Value = 12822.1578
$ws.Range("M122").Value = -507.4000000000001
$ws.Range("N122").Value = -17722.1578
$ws.Range("H123").Value = 2175.4546
$ws.Range("I123").Value = 1055
$ws.Range("J123").Value = 5163.3335
$ws.Range("K123").Value = 3165
$ws.Range("L123").Value = 15490.0005
$ws.Range("M123").Value = -715
$ws.Range("N123").Value = -20390.0005
$ws.Range("H124").Value = 530
$ws.Range("I124").Value = 530
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 1590
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 3320
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 726.87
$ws.Range("I131").Value = 615
$ws.Range("J131").Value = 729.1531
$ws.Range("K131").Value = 1845
$ws.Range("L131").Value = 2187.4593
$ws.Range("M131").Value = 3195
$ws.Range("N131").Value = -12267.4593

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2813.68
$ws.Range("I80").Value = 1665.8334
$ws.Range("K80").Value = 1665.8334
$ws.Range("M80").Value = -667.8334
$ws.Range("H83").Value = 2813.68
$ws.Range("I83").Value = 1665.8334
$ws.Range("K83").Value = 8329.166999999999
$ws.Range("M83").Value = -3337.166999999999
$ws.Range("H113").Value = 4632.0967
$ws.Range("I113").Value = 6575.3335
$ws.Range("J113").Value = 1941.4615
$ws.Range("K113").Value = 6575.3335
$ws.Range("L113").Value = 1941.4615
$ws.Range("M113").Value = -4405.3335
$ws.Range("N113").Value = -6281.461499999999
$ws.Range("H122").Value = 4142.7144
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 31280.684
$ws.Range("I132").Value = 5810
$ws.Range("K132").Value = 17430
$ws.Range("M132").Value = -14900
$ws.Range("H141").Value = 52608.547
$ws.Range("J141").Value = 52608.547
$ws.Range("L141").Value = 52608.547
$ws.Range("N141").Value = -62968.547

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2965.7368
$ws.Range("I22").Value = 4231.909
$ws.Range("K22").Value = 4231.909
$ws.Range("M22").Value = -3936.909
$ws.Range("H27").Value = 2965.7368
$ws.Range("I27").Value = 4231.909
$ws.Range("K27").Value = 4231.909
$ws.Range("M27").Value = -4124.909
$ws.Range("H46").Value = 1090.64
$ws.Range("I46").Value = 1090.5227
$ws.Range("J46").Value = 1091.5
$ws.Range("K46").Value = 1090.5227
$ws.Range("L46").Value = 1091.5
$ws.Range("M46").Value = -902.5227
$ws.Range("N46").Value = -1467.5
$ws.Range("H61").Value = 3040.6
$ws.Range("I61").Value = 1188.3
$ws.Range("J61").Value = 10449.8
$ws.Range("K61").Value = 1188.3
$ws.Range("L61").Value = 10449.8
$ws.Range("M61").Value = -986.3
$ws.Range("N61").Value = -10853.8
$ws.Range("H113").Value = 3040.6
$ws.Range("I113").Value = 1188.3
$ws.Range("J113").Value = 10449.8
$ws.Range("K113").Value = 1188.3
$ws.Range("L113").Value = 10449.8
$ws.Range("M113").Value = 981.7
$ws.Range("N113").Value = -14789.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 17500
$ws.Range("J86").Value = 17500
$ws.Range("L86").Value = 17500
$ws.Range("N86").Value = -19746
$ws.Range("H89").Value = 17500
$ws.Range("J89").Value = 17500
$ws.Range("L89").Value = 87500
$ws.Range("N89").Value = -98732
$ws.Range("H107").Value = 62500748
$ws.Range("I107").Value = 90909416
$ws.Range("J107").Value = 1671.6
$ws.Range("K107").Value = 272728248
$ws.Range("L107").Value = 5014.799999999999
$ws.Range("M107").Value = -272726328
$ws.Range("N107").Value = -8854.799999999999
$ws.Range("H132").Value = 2206.875
$ws.Range("I132").Value = 1756.8889
$ws.Range("J132").Value = 2785.4285
$ws.Range("K132").Value = 5270.6667
$ws.Range("L132").Value = 8356.2855
$ws.Range("M132").Value = -2740.6667
$ws.Range("N132").Value = -13416.2855
$ws.Range("H140").Value = 43740
$ws.Range("J140").Value = 43740
$ws.Range("L140").Value = 43740
$ws.Range("N140").Value = -54100
